$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 15: CL entry about checking the job log
$ws.Range("A15").Value = "CL"
$ws.Range("B15").Value = "Check job log"
$ws.Range("C15").Value = "WRKUSRJOB USER(BG60) STATUS(*ACTIVE)`nCan be used to view job failure point (in my case, check time deposit level up function failure point that the authority of itrtchk & itrprotyp & itrmaincl)"

# Match formatting used by the rest of the table (same wrap-text style as
# the B/C columns elsewhere, and the same row height as similarly wrapped
# rows such as row 7 and row 12).
$ws.Range("B14").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 34.5
